$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the NAICS crosswalk codes for the three rows that previously
# held compound/dotted code lists - now they reference the parent code only.
$ws.Range("B13").Value = "523"
$ws.Range("B16").Value = "53"
$ws.Range("B17").Value = "53111"

# Update the selected/active cell in the sheet view to A20.
$ws.Range("A20").Select()
